$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Feature9"
$ws.Range("A10").Value = "Feature10"

$ws.Range("F14").Select()
